$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 38/39 swap: Maker now at row 38, Hedera now at row 39
Set-TextCell "B38" "Maker"
Set-TextCell "C38" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D38" "3.173.46"
Set-TextCell "E38" "  +1.30%  "

Set-TextCell "B39" "Hedera"
Set-TextCell "C39" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D39" "0.0820"
Set-TextCell "E39" "  +6.17%  "

# Price / Volume updates for other rows
Set-TextCell "D2" "63.682.92"
Set-TextCell "E2" "  +4.42%  "
Set-TextCell "D3" "3.069.66"
Set-TextCell "E3" "  +3.29%  "
Set-TextCell "E4" "  +0.16%  "
Set-TextCell "D5" "552.06"
Set-TextCell "E5" "  +4.92%  "
Set-TextCell "D6" "139.32"
Set-TextCell "E6" "  +6.90%  "
Set-TextCell "D7" "1.00"
Set-TextCell "E7" "  +0.06%  "
Set-TextCell "D8" "3.064.34"
Set-TextCell "E8" "  +3.28%  "
Set-TextCell "D9" "0.502"
Set-TextCell "E9" "  +3.89%  "
Set-TextCell "D10" "0.150"
Set-TextCell "E10" "  +1.37%  "
Set-TextCell "D11" "6.20"
Set-TextCell "E11" "  -0.62%  "
Set-TextCell "D12" "0.456"
Set-TextCell "E12" "  +3.79%  "
Set-TextCell "D13" "0.0000227"
Set-TextCell "E13" "  +4.49%  "
Set-TextCell "D14" "35.01"
Set-TextCell "E14" "  +5.56%  "
Set-TextCell "D15" "3.576.49"
Set-TextCell "E15" "  +3.53%  "
Set-TextCell "D16" "63.789.65"
Set-TextCell "E16" "  +4.62%  "
Set-TextCell "D17" "3.080.87"
Set-TextCell "E17" "  +4.03%  "
Set-TextCell "E18" "  -0.63%  "
Set-TextCell "D19" "6.75"
Set-TextCell "E19" "  +4.72%  "
Set-TextCell "D20" "482.81"
Set-TextCell "E20" "  +5.98%  "
Set-TextCell "D21" "13.53"
Set-TextCell "E21" "  +3.43%  "
Set-TextCell "D22" "0.685"
Set-TextCell "E22" "  +2.35%  "
Set-TextCell "D23" "7.19"
Set-TextCell "E23" "  +5.97%  "
Set-TextCell "D24" "81.74"
Set-TextCell "E24" "  +5.06%  "
Set-TextCell "D25" "12.51"
Set-TextCell "E25" "  +7.00%  "
Set-TextCell "E26" "  -0.08%  "
Set-TextCell "D27" "2.76"
Set-TextCell "E27" "  +5.48%  "
Set-TextCell "D28" "7.98"
Set-TextCell "E28" "  +5.07%  "
Set-TextCell "D29" "2.00"
Set-TextCell "E29" "  +9.71%  "
Set-TextCell "D30" "1.00"
Set-TextCell "E30" "  +0.10%  "
Set-TextCell "D31" "26.06"
Set-TextCell "E31" "  +3.20%  "
Set-TextCell "D32" "1.14"
Set-TextCell "E32" "  +2.04%  "
Set-TextCell "D33" "2.43"
Set-TextCell "E33" "  +8.51%  "
Set-TextCell "D34" "5.76"
Set-TextCell "E34" "  +8.54%  "
Set-TextCell "D35" "55.67"
Set-TextCell "E35" "  +2.58%  "
Set-TextCell "D36" "5.97"
Set-TextCell "E36" "  +4.21%  "
Set-TextCell "D37" "467.79"
Set-TextCell "E37" "  +3.86%  "
Set-TextCell "D40" "0.0397"
Set-TextCell "E40" "  +4.66%  "
Set-TextCell "D41" "0.120"
Set-TextCell "E41" "  +2.19%  "
Set-TextCell "D42" "8.27"
Set-TextCell "E42" "  +4.28%  "
Set-TextCell "D43" "28.46"
Set-TextCell "E43" "  +13.95%  "
Set-TextCell "D44" "2.55"
Set-TextCell "E44" "  +7.56%  "
Set-TextCell "D45" "0.252"
Set-TextCell "E45" "  +4.26%  "
Set-TextCell "E46" "  -0.11%  "
Set-TextCell "D47" "2.03"
Set-TextCell "E47" "  +6.29%  "
Set-TextCell "E48" "  +2.69%  "
Set-TextCell "D49" "0.0₃0512"
Set-TextCell "E49" "  +2.55%  "
Set-TextCell "D50" "116.22"
Set-TextCell "E50" "  -3.27%  "
Set-TextCell "D51" "2.07"
Set-TextCell "E51" "  +6.04%  "
